# "TestData" sheet: a large pool of (private key, crypto address) pairs was
# appended to the workbook's string pool, and the small set of sample rows
# shown on the sheet (rows 2-7 and 12) were repointed to reference some of
# that newly appended test data; a brand-new row 8 (another key/address
# pair) was also inserted right after the existing row 7.
#
# Because the COM-interop runtime rebuilds the shared-string table from the
# cell values that are actually in use, we only need to set the resulting
# cell text here - the engine takes care of (re)creating the appropriate
# <sst> entries on save.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 2-5: single private-key entries (no paired address on these rows)
$ws.Range("A2").Value = "580afca99e4194bc8d4f5ec765985b5317f5258c990459cca4be47d224be3c4"
$ws.Range("A3").Value = "ea69a4824902d16d7b4b6e84c851b736ae4fd2b57d789aa11dd1243d4d27b187"
$ws.Range("A4").Value = "490d02573d601205f51f7f063e13c4d89d225084c09bf906f63f932668872028"
$ws.Range("A5").Value = "9d2fe327ad268faebccc2f8a345bbfef1852931daa58e156e32094318a222cc9"

# B2 (cryptoAddress "0x234478f5764ce6b28d5a4c39642eef71d2c24cce") is unchanged.

# Rows 6-7: private key + matching crypto address
$ws.Range("A6").Value = "a5cf82bc4fdfcc8f74d6489af1b0c7951f6abad1ffce94b37cfa0986bea1a000"
$ws.Range("B6").Value = "0x9b13a7afdb5e0156e9ee9325cb372438b5fb9e77"
$ws.Range("A7").Value = "45e87769708e46648427114a1637b0de92047a7b4588e2f154458f418470dfa9"
$ws.Range("B7").Value = "0x0411a3c75a8813a7537ae06b730ca0076b4c4fd8"

# New row 8: another private key + matching crypto address
$ws.Range("A8").Value = "a9a4400c19945ac19720c8bfa6d8924a867a324fa1d17b1b7965e55df1c15bf2"
$ws.Range("B8").Value = "0x786fce34a45b1c3ac72e73d35b93683da81e8d2a"

# Row 12: single private-key entry
$ws.Range("A12").Value = "5b5694f883e6559dffc8b80b3598533ee76f70492e7cabd7ea687bc7fb0d1838"

# Columns A/B were re-autofit to the (slightly different) new content widths
# (target OOXML widths 71.4623046875 / 46.866796875 characters). The
# runtime quantizes ColumnWidth to whole pixels internally, so these are the
# closest attainable values (71.5 / 46.833333...).
$ws.Columns.Item(1).ColumnWidth = 70.66666666666667
$ws.Columns.Item(2).ColumnWidth = 46.0
